# December 2021 - Release
# Minor updates across a few diagrams plus Matrix content refresh:
#  - bump the "last updated" date shown in A1
#  - Exchange Online Archiving moves from "Unlimited" to "1.5 TB"
#  - eDiscovery (Core) becomes available on Microsoft 365 Business Premium
#  - OneDrive for Business gains a Kiosk (2 GB) allocation on F1
#  - A5 gains Office for the Web / Office Mobile checkmarks
#  - Azure AD Premium Plan 1 checkmark removed from the standalone
#    Security add-on SKUs (F5 Security, F5 Sec+Comp, E5 Security, A5 Security)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated "last refreshed" date shown at the top of the matrix (A1) ---
# old: 2021-10-05 (serial 44474)  ->  new: 2021-12-23 (serial 44553)
$ws.Range("A1").Value = 44553

# --- Exchange Online Archiving (row 21): Unlimited -> 1.5 TB ---
$archiveCols = @("C21","D21","G21","K21","L21","M21","P21","S21","V21")
foreach ($addr in $archiveCols) {
    $ws.Range($addr).Value = "1.5 TB"
}

# --- Core eDiscovery now available on Microsoft 365 Business Premium (G19) ---
$ws.Range("G19").Value = "✔"

# --- Office for the Web / Office Mobile now checked for A5 (V50:V51) ---
$ws.Range("V50").Value = "✔"
$ws.Range("V51").Value = "✔"

# --- OneDrive for Business: F1 gets the Kiosk (2 GB) allocation (H52) ---
$ws.Range("H52").Value = "Kiosk" + [char]10 + "(2 GB)"

# --- Azure Active Directory Premium Plan 1 (rows 80-96): drop the
#     checkmark from the standalone Security SKUs (F5 Security, F5
#     Sec+Comp, E5 Security, A5 Security), keeping it in the full
#     bundles (E3, E5, A5) and Compliance SKUs untouched. ---
for ($row = 80; $row -le 96; $row++) {
    foreach ($col in @("J", "L", "N", "T")) {
        $ws.Range("$col$row").Value = $null
    }
}
